$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the newly-added log entry on row 11
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 43539
$ws.Range("D11").Value = 0.4375
$ws.Range("E11").Value = 0.57291666666666663
$ws.Range("F11").Value = 15
$ws.Range("G11").Value = 180
$ws.Range("H11").Value = "Notes"
$ws.Range("I11").Value = "Lectures Notes"

# Update the active selection to match the saved view state
$ws.Range("J17").Select()

$wb.Save()
